$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws1.Range("L5").Value = 1970.14
$ws1.Range("E25").Value = 763.17
$ws1.Range("E55").Value = "4 de 53"

$ws2.Range("F5").Value = 2691.1
$ws2.Range("F25").Value = 10627.27
$ws2.Range("F55").Value = 48705.45

$ws3.Range("D4").Value = 1012.73
$ws3.Range("E4").Value = -9.730000000000018
$ws3.Range("F4").Value = 1.009700897308076

$ws3.Range("D15").Value = 2805.36
$ws3.Range("E15").Value = 10694.64
$ws3.Range("F15").Value = 0.2078044444444445

$ws3.Range("D19").Value = 48705.45
$ws3.Range("E19").Value = 65001.00064517914
$ws3.Range("F19").Value = 0.4283437722630646

$ws3.Columns.Item(5).ColumnWidth = 23.15
